$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Simple label replacements (runs with no extra sibling markup to preserve).
# A plain Find/Replace is safe for these.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Lab1", $true, $false, $false, $false, $false, $true, 1, $false, "3.3.3", 2) | Out-Null
$d.Content.Find.Execute("Lab2", $true, $false, $false, $false, $false, $true, 1, $false, "3.4.3", 2) | Out-Null
$d.Content.Find.Execute("Lab4", $true, $false, $false, $false, $false, $true, 1, $false, "3.5.1", 2) | Out-Null

# ---------------------------------------------------------------------------
# Lab3 / Lab5: these runs also contain a <w:lastRenderedPageBreak/> marker
# right before the text. A normal Find/Replace (or any Range.Text / Delete
# operation) on that run silently drops the marker, so instead the run is
# rebuilt explicitly via InsertXML (which preserves arbitrary markup) and
# the old, now-duplicated text that InsertXML leaves behind is deleted
# afterwards (InsertXML inserts content, it does not overwrite the range).
# ---------------------------------------------------------------------------
function Replace-RunWithPageBreak([string]$oldText, [string]$newText) {
    $rng = $d.Content
    $rng.Find.Execute($oldText) | Out-Null
    $matchStart = $rng.Start
    $matchEnd = $rng.End

    $fragment = "<w:p $wNs><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr>" +
                "<w:lastRenderedPageBreak/><w:t>$newText</w:t></w:r></w:p>"

    $insertRng = $d.Range($matchStart, $matchStart)
    $insertRng.InsertXML($fragment)

    $shift = $newText.Length
    $d.Range($matchStart + $shift, $matchEnd + $shift).Delete() | Out-Null
}

Replace-RunWithPageBreak "Lab3" "3.3.4"
Replace-RunWithPageBreak "Lab5" "4.5.8"

# ---------------------------------------------------------------------------
# "Lab 6" -> "6.4.8", plus move the _GoBack bookmark so that it ends up
# right after this run (it currently sits several paragraphs later, right
# before "Lab7").
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Lab 6", $true, $false, $false, $false, $false, $true, 1, $false, "6.4.8", 2) | Out-Null

# Remove the bookmark from its old location first so we don't end up with
# two bookmarks sharing the name "_GoBack".
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Rebuild the "6.4.8" run together with the bookmark markup right after it
# (as siblings inside the same paragraph), then drop the old duplicated
# text that InsertXML left behind.
$rng = $d.Content
$rng.Find.Execute("6.4.8") | Out-Null
$matchStart = $rng.Start
$matchEnd = $rng.End

$fragment = "<w:p $wNs><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>6.4.8</w:t></w:r>" +
            '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$insertRng = $d.Range($matchStart, $matchStart)
$insertRng.InsertXML($fragment)

$shift = "6.4.8".Length
$d.Range($matchStart + $shift, $matchEnd + $shift).Delete() | Out-Null
